$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")
$ws.Activate()

# Row 18: new time-registration entry "lav SD0801 + DCD0801" / Designer, 09:00-14:00
$ws.Range("A18").Value = "lav SD0801 + DCD0801"
$ws.Range("B18").Value = "Designer"
$ws.Range("C18").Value = "2020-03-02"
$ws.Range("D18").Value = 0.375
$ws.Range("E18").Value = 0.58333333333333337

# Row 19: new time-registration entry "lav UI til UC01", 14:00-15:00
$ws.Range("A19").Value = "lav UI til UC01"
$ws.Range("C19").Value = "2020-03-02"
$ws.Range("D19").Value = 0.58333333333333337
$ws.Range("E19").Value = 0.625

# Move the active selection to E20, matching where the user left off
$ws.Range("E20").Select()
